# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-13 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 0
    4  = 5
    5  = 0
    6  = 3
    7  = 3
    8  = 1
    9  = 0
    10 = 2
    11 = 3
    12 = 3
    13 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
